$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.592.94'
$ws.Range('D3').Value = '1.878.32'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.53'
$ws.Range('E5').Value = '  +1.71%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4761'
$ws.Range('E7').Value = '  -0.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2924'
$ws.Range('E8').Value = '  +1.33%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06520'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '22.00'
$ws.Range('E10').Value = '  +2.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07735'
$ws.Range('E11').Value = '  -0.38%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.7404'
$ws.Range('E12').Value = '  +0.28%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '96.84'
$ws.Range('E13').Value = '  +0.45%  '
$ws.Range('D14').Value = '1.875.76'
$ws.Range('E14').Value = '  -0.10%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.201'
$ws.Range('E15').Value = '  +1.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '274.71'
$ws.Range('E16').Value = '  -0.51%  '
$ws.Range('D17').Value = '30.681.80'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.24'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('B19').Value = 'ShibaInu'
$ws.Range('C19').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007535'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.09%  '
$ws.Range('D21').Value = '2.122.78'
$ws.Range('E21').Value = '  -0.03%  '
$ws.Range('E22').Value = '  +0.02%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.260'
$ws.Range('E23').Value = '  +0.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.206'
$ws.Range('E24').Value = '  +0.57%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '165.53'
$ws.Range('E25').Value = '  +1.19%  '
$ws.Range('B26').Value = 'Cosmos'
$ws.Range('C26').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.195'
$ws.Range('E26').Value = '  -0.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.89'
$ws.Range('E27').Value = '  -0.36%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.916'
$ws.Range('E28').Value = '  -1.99%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.09865'
$ws.Range('E29').Value = '  -0.91%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.340'
$ws.Range('E30').Value = '  -2.21%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.503'
$ws.Range('E31').Value = '  -0.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.285'
$ws.Range('E32').Value = '  -0.71%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.116'
$ws.Range('E33').Value = '  +0.81%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.04818'
$ws.Range('E34').Value = '  +1.49%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.126'
$ws.Range('E35').Value = '  +0.19%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6966'
$ws.Range('E36').Value = '  -0.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.717'
$ws.Range('E37').Value = '  +0.03%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01869'
$ws.Range('E38').Value = '  +0.52%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.765'
$ws.Range('E39').Value = '  +0.40%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '6.277'
$ws.Range('E40').Value = '  -0.06%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '73.55'
$ws.Range('E41').Value = '  +5.87%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.996'
$ws.Range('E42').Value = '  +4.58%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4235'
$ws.Range('E43').Value = '  +1.43%  '
$ws.Range('E44').Value = '  +0.09%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.8380'
$ws.Range('E45').Value = '  -0.39%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '102.08'
$ws.Range('E46').Value = '  +0.33%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.327'
$ws.Range('E47').Value = '  +0.66%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.022'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '35.43'
$ws.Range('E49').Value = '  +0.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '912.61'
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05686'
